# Auto-update draw results: append the 2025-11-02 Pick 3 draw as row 47.
#
# All five columns in this sheet are stored as literal text (Date, Game,
# Phase, Result, InsertedAt), matching every prior row. Two of the new
# values (the ISO date "2025-11-02" and the all-digit phase code
# "251102") look like a date / number to Excel's input parser, so a
# plain .Value assignment would silently convert them to a real date
# serial / numeric value. Prefixing with a leading apostrophe is the
# standard Excel technique to force those two values to stay literal
# text, exactly like the existing rows above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 47

$ws.Cells.Item($row, 1).Value = "'2025-11-02"
$ws.Cells.Item($row, 2).Value = "Pick 3"
$ws.Cells.Item($row, 3).Value = "'251102"
$ws.Cells.Item($row, 4).Value = "1-1-9"
$ws.Cells.Item($row, 5).Value = "2025-11-02T21:36:06.525+04:00"
